$d = $word.ActiveDocument

# The "SMARTREWARDS FAQ'S" heading paragraph is kept, but the whole FAQ
# block that used to follow it (all the Q&A paragraphs, down to the final
# answer about when SmartReward payments go out) is removed, leaving the
# heading immediately followed by the final (empty) paragraph of the doc.

# Locate the first paragraph of the block to remove: the first FAQ
# question, right after the "SMARTREWARDS FAQ'S" heading.
$startRange = $d.Content.Duplicate
$startFound = $startRange.Find.Execute(
    "How much SmartRewards can I get each month?",
    $true, $true, $false, $false, $false, $true, 1, $false, "", 0)

# Locate the last paragraph of the block to remove: the final FAQ answer
# about SmartReward payments going out.
$endRange = $d.Content.Duplicate
$endFound = $endRange.Find.Execute(
    "every second block 1000 addresses will get paid.",
    $true, $true, $false, $false, $false, $true, 1, $false, "", 0)

if ($startFound -and $endFound) {
    # Expand to the full paragraphs (including their paragraph marks) so
    # the entire block -- and nothing else -- is deleted.
    $delStart = $startRange.Paragraphs(1).Range.Start
    $delEnd = $endRange.Paragraphs(1).Range.End

    $toDelete = $d.Range($delStart, $delEnd)
    $toDelete.Delete()
}
